$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.878.19'
$ws.Range("E2").Value = '  -0.50%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.859.38'

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.48'
$ws.Range("E5").Value = '  -0.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9998'

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5026'
$ws.Range("E7").Value = '  -1.62%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3645'
$ws.Range("E8").Value = '  -2.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07170'
$ws.Range("E9").Value = '  +0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8922'
$ws.Range("E10").Value = '  +0.19%  '

$ws.Range("E11").Value = '  +0.44%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.876.22'
$ws.Range("E12").Value = '  +1.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07484'
$ws.Range("E13").Value = '  -0.74%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '93.55'
$ws.Range("E14").Value = '  +5.49%  '

$ws.Range("E16").Value = '  -0.05%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008490'
$ws.Range("E17").Value = '  +1.25%  '

$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.928.34'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("E21").Value = '  -0.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.113.56'
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("E23").Value = '  -1.46%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.410'
$ws.Range("E24").Value = '  -0.99%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.54'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.782'
$ws.Range("E26").Value = '  -3.28%  '

$ws.Range("E27").Value = '  -0.60%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.081'
$ws.Range("E28").Value = '  -0.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.02'
$ws.Range("E29").Value = '  +0.17%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.683'
$ws.Range("E30").Value = '  +0.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.666'
$ws.Range("E31").Value = '  +0.34%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09212'
$ws.Range("E32").Value = '  +1.88%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05136'
$ws.Range("E33").Value = '  +0.40%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7462'
$ws.Range("E34").Value = '  +2.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.953'
$ws.Range("E35").Value = '  -3.86%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.150'
$ws.Range("E36").Value = '  -0.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.251'
$ws.Range("E37").Value = '  +6.63%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.564'
$ws.Range("E38").Value = '  +2.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02002'
$ws.Range("E39").Value = '  -2.25%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5552'
$ws.Range("E40").Value = '  +4.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.069'
$ws.Range("E41").Value = '  -0.23%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.539'
$ws.Range("E42").Value = '  -0.77%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.27'
$ws.Range("E43").Value = '  +1.97%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.507'
$ws.Range("E44").Value = '  +2.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4670'
$ws.Range("E46").Value = '  +1.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9995'
$ws.Range("E47").Value = '  -0.17%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.01'
$ws.Range("E48").Value = '  -0.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.558'
$ws.Range("E49").Value = '  -0.41%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.64'
$ws.Range("E50").Value = '  -0.15%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '62.91'
$ws.Range("E51").Value = '  -1.90%  '
